$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the paragraph "whose, possessive 's - Whose is this? ..."
#    New "at, in, on" paragraphs are inserted right after it.
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*whose, possessive*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the 'whose, possessive' paragraph"
}

$whoseParagraph = $d.Paragraphs.Item($targetIndex)

# Insert the "prepositions of time" paragraph right after it.
$whoseParagraph.Range.InsertParagraphAfter()
$timeParagraph = $d.Paragraphs.Item($targetIndex + 1)
$timeTextRange = $d.Range($timeParagraph.Range.Start, $timeParagraph.Range.End - 1)
$timeTextRange.Text = "at, in, on – prepositions of time"

# Insert the "prepositions of place" paragraph right after that one.
$timeParagraph.Range.InsertParagraphAfter()
$placeParagraph = $d.Paragraphs.Item($targetIndex + 2)
$placeTextRange = $d.Range($placeParagraph.Range.Start, $placeParagraph.Range.End - 1)
$placeTextRange.Text = "at, in, on – prepositions of place"

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the "whose, ..."
#    paragraph to just after the page-break run in the paragraph
#    that now follows the (still) empty paragraph.
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# Find the paragraph that contains the page break (w:br w:type="page"),
# which is two paragraphs after the empty paragraph that used to sit
# right after "whose, possessive ...".
$pageBreakIndex = -1
for ($i = $targetIndex; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$([char]12)*") {
        $pageBreakIndex = $i
        break
    }
}
if ($pageBreakIndex -eq -1) {
    throw "Could not find the page-break paragraph"
}

$pageBreakParagraph = $d.Paragraphs.Item($pageBreakIndex)
$targetPos = $pageBreakParagraph.Range.End - 1

# Collapsed (zero-length) ranges placed exactly at a run/paragraph-mark
# boundary confuse Bookmarks.Add in this host, so work around it: insert
# a throw-away character, wrap it with the bookmark, then delete the
# character again. The bookmark collapses correctly in its place.
$insertionPoint = $d.Range($targetPos, $targetPos)
$insertionPoint.InsertAfter("X")
$placeholderRange = $d.Range($targetPos, $targetPos + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$placeholderRange.Text = ""
